$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.969.33"
$ws.Range("E2").Value = "  -2.72%  "

# Row 3
$ws.Range("D3").Value = "1.792.71"
$ws.Range("E3").Value = "  -3.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.36%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.07"
$ws.Range("E5").Value = "  -2.32%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4168"
$ws.Range("E7").Value = "  -3.30%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3547"
$ws.Range("E8").Value = "  -4.36%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07049"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8393"
$ws.Range("E10").Value = "  -4.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.11"
$ws.Range("E11").Value = "  -4.29%  "

# Row 12
$ws.Range("D12").Value = "1.790.93"
$ws.Range("E12").Value = "  -3.70%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.278"
$ws.Range("E13").Value = "  -3.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.324"
$ws.Range("E14").Value = "  -4.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06734"
$ws.Range("E15").Value = "  -3.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.47"
$ws.Range("E17").Value = "  -2.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008679"
$ws.Range("E18").Value = "  -4.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.59%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.01"
$ws.Range("E20").Value = "  -3.65%  "

# Row 21
$ws.Range("D21").Value = "27.193.50"
$ws.Range("E21").Value = "  -2.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.041"
$ws.Range("E22").Value = "  -0.98%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").Value = "2.085.39"
$ws.Range("E24").Value = "  -0.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.940"
$ws.Range("E25").Value = "  -1.32%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.67"
$ws.Range("E26").Value = "  -1.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.07"
$ws.Range("E27").Value = "  -2.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.983"
$ws.Range("E28").Value = "  -6.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.65"
$ws.Range("E29").Value = "  -2.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.631"
$ws.Range("E30").Value = "  -11.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08896"
$ws.Range("E31").Value = "  -0.38%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7124"
$ws.Range("E32").Value = "  -9.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.848"
$ws.Range("E33").Value = "  -4.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.294"
$ws.Range("E34").Value = "  -6.78%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.068"
$ws.Range("E36").Value = "  -8.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.074"
$ws.Range("E37").Value = "  -3.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01894"
$ws.Range("E38").Value = "  -3.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05095"
$ws.Range("E39").Value = "  -6.10%  "

# Row 40
$ws.Range("E40").Value = "  -4.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4924"
$ws.Range("E41").Value = "  -5.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.575"
$ws.Range("E42").Value = "  -9.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.010"
$ws.Range("E43").Value = "  -11.13%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.987"
$ws.Range("E44").Value = "  -7.83%  "

# Row 45
$ws.Range("E45").Value = "  +0.49%  "

# Row 46
$ws.Range("E46").Value = "  -4.59%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.08"
$ws.Range("E47").Value = "  -2.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06295"
$ws.Range("E48").Value = "  -4.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4489"
$ws.Range("E49").Value = "  -6.31%  "

# Row 50
$ws.Range("E50").Value = "  -4.94%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.67"
$ws.Range("E51").Value = "  -5.04%  "
